$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.874.72'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('D3').Value = '3.385.95'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '564.05'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.57'
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('E7').Value = '  +0.74%  '
$ws.Range('D8').Value = '3.381.60'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +2.61%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.94'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = '3.930.64'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').Value = '3.382.36'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').Value = '65.711.23'
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.90'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '464.43'
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.92'
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.82'
$ws.Range('E24').Value = '  +10.05%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.11'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.48'
$ws.Range('E26').Value = '  +2.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.93'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.67'
$ws.Range('E28').Value = '  -1.87%  '
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.13'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.62'
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '580.86'
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '62.41'
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +1.89%  '
$ws.Range('E38').Value = '  +1.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.07'
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.379'
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('D42').Value = '3.104.07'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.18'
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.45'
$ws.Range('E47').Value = '  -0.87%  '
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.65'
$ws.Range('E49').Value = '  +1.22%  '
$ws.Range('E50').Value = '  +9.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.17'
$ws.Range('E51').Value = '  +9.79%  '
